$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-22 19:10:01"
$zhcn.Range("H2").Value = "2016-03-22 19:10:40"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-22 19:10:08"
$dede.Range("H2").Value = "2016-03-22 19:10:48"
